$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -22.192
$ws.Range("C4").Value = -12.925

$ws.Range("C5").Value = -12.413

$ws.Range("A7").Value = -20.987

$ws.Range("C8").Value = -12.8

$ws.Range("A16").Value = -21.355
$ws.Range("C16").Value = -12.346
